$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dutch")

# Insert a new row before row 87, shifting existing rows 87-101 down.
$ws.Rows.Item(87).Insert()

# Fill the new row 87 with the new minimal-stress-pair entry.
$ws.Cells.Item(87, 1).Value = 4
$ws.Cells.Item(87, 2).Value = "achterwegen"
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(87, 4).Value = "back roads"
$ws.Cells.Item(87, 5).Value = "."
$ws.Cells.Item(87, 6).Value = "achterwege"
$ws.Cells.Item(87, 7).Value = 3
$ws.Cells.Item(87, 8).Value = "omitted/left out (adv.)"

# Update the view so the newly added row is visible/selected.
$excel.ActiveWindow.ScrollRow = 78
$ws.Range("H87").Select()
